$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.455.20'
$ws.Range("E2").Value = '  -2.31%  '

$ws.Range("D3").Value = '3.149.52'
$ws.Range("E3").Value = '  -3.90%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "'526.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.97%  '

$ws.Range("D6").Value = "'135.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.34%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '3.146.88'
$ws.Range("E8").Value = '  -3.91%  '

$ws.Range("E9").Value = '  -4.18%  '

$ws.Range("D10").Value = "'7.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.55%  '

$ws.Range("D11").Value = "'0.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.31%  '

$ws.Range("D12").Value = "'0.377"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.40%  '

$ws.Range("D13").Value = '3.681.00'
$ws.Range("E13").Value = '  -3.98%  '

$ws.Range("E14").Value = '  -0.88%  '

$ws.Range("D15").Value = "'25.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.22%  '

$ws.Range("D16").Value = '3.142.00'
$ws.Range("E16").Value = '  -3.77%  '

$ws.Range("D17").Value = '58.395.94'
$ws.Range("E17").Value = '  -2.49%  '

$ws.Range("E18").Value = '  -6.32%  '

$ws.Range("E19").Value = '  -4.76%  '

$ws.Range("D20").Value = "'13.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.37%  '

$ws.Range("D21").Value = "'7.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.38%  '

$ws.Range("D22").Value = "'343.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.45%  '

$ws.Range("E23").Value = '  -0.14%  '

$ws.Range("D24").Value = "'0.510"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.50%  '

$ws.Range("D25").Value = "'67.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.00%  '

$ws.Range("D26").Value = '3.269.81'
$ws.Range("E26").Value = '  -3.88%  '

$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("D28").Value = '0.0₃0953'
$ws.Range("E28").Value = '  -6.56%  '

$ws.Range("E29").Value = '  +0.13%  '

$ws.Range("D30").Value = "'6.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.68%  '

$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = "'1.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.97%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = "'6.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.60%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = "'1.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.55%  '

$ws.Range("D35").Value = "'21.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.29%  '

$ws.Range("E36").Value = '  -3.69%  '

$ws.Range("D37").Value = "'157.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.81%  '

$ws.Range("D38").Value = "'6.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.99%  '

$ws.Range("D39").Value = "'1.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.40%  '

$ws.Range("E40").Value = '  -4.89%  '

$ws.Range("D41").Value = '3.173.83'
$ws.Range("E41").Value = '  -3.95%  '

$ws.Range("D42").Value = "'40.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.14%  '

$ws.Range("D43").Value = "'23.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.10%  '

$ws.Range("D44").Value = "'1.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.53%  '

$ws.Range("D45").Value = "'0.694"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.65%  '

$ws.Range("E46").Value = '  -4.11%  '

$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("D48").Value = "'1.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.49%  '

$ws.Range("D49").Value = '2.280.39'
$ws.Range("E49").Value = '  -1.44%  '

$ws.Range("D50").Value = "'6.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.67%  '

$ws.Range("D51").Value = "'20.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.19%  '
